# Fix article codes that were mistakenly typed with a trailing "p"
# (commit: "Izbrisani p pri sifrah - NAROBE VPISALI")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# MAG167p .. MAG177p  ->  MAG167 .. MAG177
$ws.Range("B128").Value = "MAG167"
$ws.Range("B129").Value = "MAG168"
$ws.Range("B130").Value = "MAG169"
$ws.Range("B131").Value = "MAG170"
$ws.Range("B132").Value = "MAG171"
$ws.Range("B133").Value = "MAG172"
$ws.Range("B134").Value = "MAG173"
$ws.Range("B135").Value = "MAG174"
$ws.Range("B136").Value = "MAG175"
$ws.Range("B137").Value = "MAG176"
$ws.Range("B138").Value = "MAG177"

# MAG187p .. MAG191p  ->  MAG187 .. MAG191
$ws.Range("B148").Value = "MAG187"
$ws.Range("B149").Value = "MAG188"
$ws.Range("B150").Value = "MAG189"
$ws.Range("B151").Value = "MAG190"
$ws.Range("B152").Value = "MAG191"

# MAG135p .. MAG152p  ->  MAG135 .. MAG152
$ws.Range("B98").Value = "MAG135"
$ws.Range("B99").Value = "MAG136"
$ws.Range("B100").Value = "MAG137"
$ws.Range("B101").Value = "MAG138"
$ws.Range("B102").Value = "MAG139"
$ws.Range("B103").Value = "MAG140"
$ws.Range("B104").Value = "MAG141"
$ws.Range("B105").Value = "MAG142"
$ws.Range("B106").Value = "MAG143"
$ws.Range("B107").Value = "MAG144"
$ws.Range("B108").Value = "MAG145"
$ws.Range("B109").Value = "MAG146"
$ws.Range("B110").Value = "MAG147"
$ws.Range("B111").Value = "MAG148"
$ws.Range("B112").Value = "MAG149"
$ws.Range("B113").Value = "MAG150"
$ws.Range("B114").Value = "MAG151"
$ws.Range("B115").Value = "MAG152"

$ws.Range("B115").Select()
